$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-19.01053083759864, 1.954084305947056, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(-19.01053083759864, 1.966719035032424, 1.704034158973791, -19.01053083759864, 3.426375098842972, -19.01053083759864, 1.521974434756883, -19.01053083759864, 0.4894026154519405, -19.01053083759864),
    @(-19.01053083759864, 1.645846014655267, -19.01053083759864, -19.01053083759864, -19.01053083759864, 2.819308396477175, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(2.454940409745085, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, 1.779958460916872, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(3.859677519232898, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, 2.273774226407448),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, 2.929510636850866, -19.01053083759864, 2.764081307399505, -19.01053083759864, -19.01053083759864, -19.01053083759864, 1.957548802016823),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, 2.515242764311107, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, 1.728574834053736, 1.768955287871083),
    @(-19.01053083759864, -19.01053083759864, 1.495518094356097, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, 1.939430487039074),
    @(-19.01053083759864, -19.01053083759864, 1.754033843350523, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, 2.03492555574605, -19.01053083759864),
    @(-19.01053083759864, 2.190993070124682, 1.8408895797645, -19.01053083759864, -19.01053083759864, -19.01053083759864, 2.001423201394328, -19.01053083759864, 2.547337148032323, -19.01053083759864),
    @(-19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, -19.01053083759864, 1.816325615776167, -19.01053083759864, 2.416446884963039, -19.01053083759864),
    @(-19.01053083759864, -19.01053083759864, 2.054088924284129, -19.01053083759864, -19.01053083759864, -19.01053083759864, 1.709323357390094, 4.321925499760941, -19.01053083759864, -19.01053083759864),
    @(-19.01053083759864, 1.059905923966017, 1.492666976870667, -19.01053083759864, 3.209319800997702, -19.01053083759864, 1.66664596416075, -19.01053083759864, -19.01053083759864, 2.013096967926595),
    @(-19.01053083759864, 1.28462487312774, -19.01053083759864, 1.691690050003025, -19.01053083759864, 2.62025219260526, 1.65879457137238, -19.01053083759864, -19.01053083759864, -19.01053083759864),
)

for ($r = 0; $r -lt $values.Count; $r++) {
    for ($c = 0; $c -lt $values[$r].Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $values[$r][$c]
    }
}
